# Update the "取得日時" (acquired timestamp) column for the appended rows
# from 2026-01-26 01:45:08 to 2026-01-26 02:11:02 on the "ランサーズ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2026-01-26 01:45:08") {
        $cell.Value2 = "2026-01-26 02:11:02"
    }
}
